$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 3545.742
$ws.Range("I113").Value = 2970.6155
$ws.Range("J113").Value = 3961.111
$ws.Range("K113").Value = 2970.6155
$ws.Range("L113").Value = 3961.111
$ws.Range("M113").Value = 283.3845000000001
$ws.Range("N113").Value = -10469.111
$ws.Range("H127").Value = 994.35297
$ws.Range("I127").Value = 661.2857
$ws.Range("J127").Value = 1227.5
$ws.Range("K127").Value = 1983.8571
$ws.Range("L127").Value = 3682.5
$ws.Range("M127").Value = 2976.1429
$ws.Range("N127").Value = -13602.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6303.73
$ws.Range("I32").Value = 3285.75
$ws.Range("J32").Value = 22148.125
$ws.Range("K32").Value = 3285.75
$ws.Range("L32").Value = 22148.125
$ws.Range("M32").Value = -2998.75
$ws.Range("N32").Value = -22722.125
$ws.Range("H61").Value = 2832.6482
$ws.Range("I61").Value = 1706.1842
$ws.Range("J61").Value = 5508
$ws.Range("K61").Value = 1706.1842
$ws.Range("L61").Value = 5508
$ws.Range("M61").Value = -1494.1842
$ws.Range("N61").Value = -5932
$ws.Range("H136").Value = 2832.6482
$ws.Range("I136").Value = 1706.1842
$ws.Range("J136").Value = 5508
$ws.Range("K136").Value = 5118.5526
$ws.Range("L136").Value = 16524
$ws.Range("M136").Value = -2568.5526
$ws.Range("N136").Value = -21624

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 364.83334
$ws.Range("I64").Value = 115.57143
$ws.Range("K64").Value = 115.57143
$ws.Range("M64").Value = 109.42857
$ws.Range("H67").Value = 364.83334
$ws.Range("I67").Value = 115.57143
$ws.Range("K67").Value = 115.57143
$ws.Range("M67").Value = 664.42857
$ws.Range("H94").Value = 974.775
$ws.Range("I94").Value = 869.76
$ws.Range("J94").Value = 1149.8
$ws.Range("K94").Value = 869.76
$ws.Range("L94").Value = 1149.8
$ws.Range("M94").Value = -418.76
$ws.Range("N94").Value = -2051.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2466.6667
$ws.Range("I16").Value = 2780
$ws.Range("K16").Value = 2780
$ws.Range("M16").Value = -2493
$ws.Range("H99").Value = 86662.25
$ws.Range("I99").Value = 102994.7
$ws.Range("J99").Value = 5000
$ws.Range("K99").Value = 102994.7
$ws.Range("L99").Value = 5000
$ws.Range("M99").Value = -101496.7
$ws.Range("N99").Value = -7996
$ws.Range("H105").Value = 617.11536
$ws.Range("I105").Value = 602.1739
$ws.Range("J105").Value = 731.6667
$ws.Range("K105").Value = 602.1739
$ws.Range("L105").Value = 731.6667
$ws.Range("M105").Value = 1144.8261
$ws.Range("N105").Value = -4225.6667
$ws.Range("H113").Value = 2466.6667
$ws.Range("I113").Value = 2780
$ws.Range("K113").Value = 2780
$ws.Range("M113").Value = -610
$ws.Range("H126").Value = 86662.25
$ws.Range("I126").Value = 102994.7
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 308984.1
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -306514.1
$ws.Range("N126").Value = -19940

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 5742.857
$ws.Range("I56").Value = 5742.857
$ws.Range("K56").Value = 5742.857
$ws.Range("M56").Value = -5212.857
$ws.Range("H92").Value = 912.9231
$ws.Range("J92").Value = 924
$ws.Range("L92").Value = 2772
$ws.Range("N92").Value = -5268
$ws.Range("H113").Value = 596.7442
$ws.Range("I113").Value = 585.8570999999999
$ws.Range("J113").Value = 602
$ws.Range("K113").Value = 1757.5713
$ws.Range("L113").Value = 1806
$ws.Range("M113").Value = 412.4287000000002
$ws.Range("N113").Value = -6146
$ws.Range("H116").Value = 2603.2
$ws.Range("I116").Value = 604
$ws.Range("J116").Value = 4602.4
$ws.Range("K116").Value = 1812
$ws.Range("L116").Value = 13807.2
$ws.Range("M116").Value = 1630
$ws.Range("N116").Value = -20691.2
$ws.Range("H131").Value = 1791.5454
$ws.Range("J131").Value = 1478.1404
$ws.Range("L131").Value = 4434.4212
$ws.Range("N131").Value = -14514.4212

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 858.1429000000001
$ws.Range("I13").Value = 500.25
$ws.Range("J13").Value = 1335.3334
$ws.Range("K13").Value = 500.25
$ws.Range("L13").Value = 1335.3334
$ws.Range("M13").Value = -361.25
$ws.Range("N13").Value = -1613.3334
$ws.Range("H102").Value = 2762.3948
$ws.Range("I102").Value = 2925.2058
$ws.Range("J102").Value = 1378.5
$ws.Range("K102").Value = 2925.2058
$ws.Range("L102").Value = 1378.5
$ws.Range("M102").Value = -1303.2058
$ws.Range("N102").Value = -4622.5
$ws.Range("H122").Value = 1102.9
$ws.Range("I122").Value = 1159.1666
$ws.Range("J122").Value = 1018.5
$ws.Range("K122").Value = 3477.4998
$ws.Range("L122").Value = 3055.5
$ws.Range("M122").Value = -1027.4998
$ws.Range("N122").Value = -7955.5
$ws.Range("H126").Value = 2227.7273
$ws.Range("I126").Value = 2132.6875
$ws.Range("J126").Value = 2481.1667
$ws.Range("K126").Value = 6398.0625
$ws.Range("L126").Value = 7443.500100000001
$ws.Range("M126").Value = -3928.0625
$ws.Range("N126").Value = -12383.5001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2150.2307
$ws.Range("I7").Value = 2053.4
$ws.Range("J7").Value = 2282.2727
$ws.Range("K7").Value = 2053.4
$ws.Range("L7").Value = 2282.2727
$ws.Range("M7").Value = -1941.4
$ws.Range("N7").Value = -2506.2727
$ws.Range("H40").Value = 2887.76
$ws.Range("I40").Value = 2806.4443
$ws.Range("J40").Value = 3096.8572
$ws.Range("K40").Value = 2806.4443
$ws.Range("L40").Value = 3096.8572
$ws.Range("M40").Value = -2670.4443
$ws.Range("N40").Value = -3368.8572
$ws.Range("H126").Value = 2150.2307
$ws.Range("I126").Value = 2053.4
$ws.Range("J126").Value = 2282.2727
$ws.Range("K126").Value = 6160.200000000001
$ws.Range("L126").Value = 6846.8181
$ws.Range("M126").Value = -3690.200000000001
$ws.Range("N126").Value = -11786.8181

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2203.875
$ws.Range("I126").Value = 1166.6666
$ws.Range("J126").Value = 2826.2
$ws.Range("K126").Value = 3499.9998
$ws.Range("L126").Value = 8478.599999999999
$ws.Range("M126").Value = -1029.9998
$ws.Range("N126").Value = -13418.6
$ws.Range("H136").Value = 12988392
$ws.Range("J136").Value = 2125.111
$ws.Range("L136").Value = 6375.333
$ws.Range("M136").Value = -6375.333
$ws.Range("N136").Value = -11475.333
